$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# New Week_Start_Date (col B) and MyForecast (col D) values for rows 2..17.
$weekStartDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(215, 216, 217, 223, 137, 138, 131, 133, 130, 130, 122, 128, 124, 119, 120, 119)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    # Force text so "YYYY-MM-DD" stays a literal string instead of being
    # auto-parsed into a date serial number (matches the source data, which
    # stores these as plain text, not real dates).
    $cellB = $ws1.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $weekStartDates[$i]
    $ws1.Cells.Item($row, 4).Value = $myForecast[$i]
}

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

# Every value in column B on this sheet is stored as plain text in the
# source data (even the purely-numeric-looking ones, e.g. "159"), so force
# text formatting before assigning new values to avoid Excel auto-converting
# them into numbers or dates.
$summaryRows = @(2, 5, 6, 8, 9, 10, 11, 12, 13, 14, 15)
foreach ($r in $summaryRows) {
    $ws2.Cells.Item($r, 2).NumberFormat = "@"
}

$ws2.Range("B2").Value = "2022-12-25 to 2025-01-05"
$ws2.Range("B5").Value = "158"
$ws2.Range("B6").Value = "139"
$ws2.Range("B8").Value = "16779 units"
$ws2.Range("B9").Value = "2401"
$ws2.Range("B10").Value = "1410"
$ws2.Range("B11").Value = "871"
$ws2.Range("B12").Value = "223"
$ws2.Range("B13").Value = "2025-02-02"
$ws2.Range("B14").Value = "119"
$ws2.Range("B15").Value = "2025-04-27"
